$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For every data row (2 through 70), the "birth_year" column (Q) is
# decremented by 1 and the "age_y" column (S) is incremented by 1.
for ($row = 2; $row -le 70; $row++) {
    $qCell = $ws.Range("Q$row")
    $qCell.Value = $qCell.Value() - 1

    $sCell = $ws.Range("S$row")
    $sCell.Value = $sCell.Value() + 1
}
